$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Ativo" column (E) used to store the boolean-like text VERDADEIRO/FALSO.
# Replace it with the clearer ATIVO/INATIVO labels (design/export bug fix).
$ws.Range("E2").Value = "ATIVO"
$ws.Range("E3").Value = "INATIVO"
$ws.Range("E4").Value = "ATIVO"
$ws.Range("E5").Value = "ATIVO"
$ws.Range("E6").Value = "INATIVO"
$ws.Range("E7").Value = "ATIVO"
$ws.Range("E8").Value = "ATIVO"
